$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 4.62 = 18418.01 pesos", "1000 Bs = 4.56 = 18129.23 pesos")
$text = $text.Replace("18418.01 pesos = 4.61 = 964.27 Bs", "18129.23 pesos = 4.53 = 954.53 Bs")
$cell.Value = $text

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 219.37
$wsTasas.Range("O10").Value = 3977.01
$wsTasas.Range("N12").Value = 3998
$wsTasas.Range("O12").Value = 210.5
